# Agregado soporte para Brasil (QUILOGRAMA LIQUIDO) y cálculo de Toneladas Finales
#
# For rows 100-137 (Brasil imports reported in QUILOGRAMA LIQUIDO), normalize
# the unit of measure to KILOGRAMOS and compute "Toneladas Finales" (Y) from
# "Cantidad Comercial" (W, in kilograms) by dividing by 1000.
# A subset of these rows additionally has "Aplica?" (A) flipped from NO to SI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Aplica?" flag (column A) changes from "NO" to "SI"
$siRows = @(100, 102, 105, 106, 107, 109, 110, 115, 116, 121, 123, 125, 126, 127, 128, 130, 131, 132, 133, 135, 137)

for ($r = 100; $r -le 137; $r++) {
    if ($siRows -contains $r) {
        $ws.Cells.Item($r, 1).Value = "SI"
    }

    $cantidad = $ws.Cells.Item($r, 23).Value2   # W: Cantidad Comercial (kg)

    $ws.Cells.Item($r, 24).Value = "KILOGRAMOS"  # X: Unidad de Medida
    $ws.Cells.Item($r, 25).Value = $cantidad / 1000  # Y: Toneladas Finales
}
